$d = $word.ActiveDocument

# 1. Update the date in the first line
$d.Content.Find.Execute("⚡️🚀המאמר היומי של מייק 01.08.24: ⚡️🚀", $true, $false, $false, $false, $false, $true, 1, $false, "⚡️🚀המאמר היומי של מייק 31.07.24: ⚡️🚀", 2) | Out-Null

# 2. Update the paper title and add a manual line break after it (char 11 = vertical tab = w:br)
$newTitle = "DENOISING DIFFUSION IMPLICIT MODELS" + [char]11
$d.Content.Find.Execute("IMPROVED TECHNIQUES FOR TRAINING CONSISTENCY MODELS", $true, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2) | Out-Null

# 3. Rewrite the intro paragraph
$d.Content.Find.Execute("היום סוקרים קצרות עוד מאמר בנושא קרוב לליבי - המשך של המאמר שסקרנו לפני בערך שבוע הנקרא ״consistency models`". אם אתם זוכרים מודל קונסיסטנטי הוא שייך למשפחת מודלי דיפוזיה (כלומר הוא מתואר על ידי משוואת הדיפוזיה). אחת הבעיות של מודלי דיפוזיה קלאסיים (כמו DDPM) היא איטיות של גנרוט דאטה. הדאטה נוצר באמצעות תהליך denoising הדרגתי - מתחילים עם רעש גאוסי ומסירים אותו לאט לאט.", $true, $false, $false, $false, $false, $true, 1, $false, "זה מאמר לא חדש (אוקטובר 2022) אך חשוב מאוד בתחום של מודלי דיפוזיה. מאמר עם רעיון מאוד אלגנטי המלווה במתמטיקה די רצינית. אנסה לסקור אותו קצרות כי כאמור יש בו עומק מתמטי לא קטן אך עדיין ניתן להעביר את הרעיון העיקרי בלי לצלול יותר מדי לעומק.", 2) | Out-Null

# 4. Rewrite the paragraph about the forward/backward diffusion process
$d.Content.Find.Execute("כדי להתמודד עם הבעיה הזו הוצעו כמה שיטות ואחת מהן DDIM סקרנו אתמול. השנייה היא מודלים קונסיסטנטיים(CM) שניתן להגדיר אותם כי מודל שונה (אך דומה) ממודל דיפוזיה קלאסי. בעיקרון ב-CM אנו מאמנים מודל להסיר רעש מכל פיסת דאטה מורעש באיטרציה t כך שהתוצאה תמיד תהיה פיסת הדאטה מקורית (ללא רעש). מכאן בא שם של המודל: קונסיסטנטי. ", $true, $false, $false, $false, $false, $true, 1, $false, "כמו שאתם זוכרים במודלי דיפוזיה גנרטיביים יש לנו שני תהליכים: הקדמי והאחורי. תהליך הקדמי הוא הרעשה הדרגתית של דאטה והתהליך האחורי הוא הורדה הדרגתית של הרעש מהדאטה באמצעות מודל שאומן לצורך זה על דאטהסט מסוים. למעשה מודל כזה מאפשר ליצור דאטה מרעש טהור בצורה הדרגתית. הבעיה בתהליך הזה כמובן זה הזמן שזה לוקח כי צריך די הרבה איטרציות של denoising כדי להגיע מרעש לדאטה איכותי.", 2) | Out-Null

# 5. Rewrite the paragraph about reducing iterations
$d.Content.Find.Execute("איך זה למעשה נעשה? יש שתי דרכים עיקריות לאמן CM. דרך אחת מסתמכת על מודל המשערך את מה שנקרא score function שהיא לוגריתם של פונקציית ההסתברות של הדאטה המורעש באיטרציה t. ידוע כי תהליך גנרוט של דאטה במודלי דיפוזיה (כלומר denoising) מתואר על ידי משוואת זרימה (דיפרנציאלית) שמתאר את המסלול של דאטה מהרעש עד הדאטה הנקי. ו- score function מופיע במשוואת זרימה זו. אז השיטה הראשונה ממזערת את המרחק בין שערוך x_0 (הדאטה הנקי) מ x_t+1 לבין שערוך של x_0 מ- x_t כאשר x_t מחושב ממשוואת הזרימה (איטרציה אחת של אוילר של משוואת הזרימה שכבר הזכרנו). ו", $true, $false, $false, $false, $false, $true, 1, $false, "המאמר מציע דרך להקטין את מספר האיטרציות בדרך די מקורית. כמו שאתם זוכרים תהליך ההרעשה (הקדמי) במודלי דיפוזיה רגילים הוא מרקובי, כלומר הדאטה באיטרציה t מוגדר (מבחינת התפלגות) על ידי הדאטה המורעש מאיטרציה t-1 בלבד כל. המאמר הורס את ההנחה הזו ומגדיר תהליך קידמי לא מרקוב כאשר הדאטה באיטרציה t מוגדר לא רק על ידי הדאטה באיטרציה t-1 אלא גם על ידי הדאטה הנקי (x_0). ", 2) | Out-Null

# 6. Rewrite the paragraph about the deterministic process
$d.Content.Find.Execute("דרך אגב שערוך של score function די שקול לשערוך של הרעש הנוסף (לדאטה) במודלי הדיפוזיה הסטנדרטיים. הדרך השנייה ״ליצור״ את x_t+1 היא לשערך את x_0 מ-x_t+1 ולהוסיף רעש (כמו באיטרציה t). ", $true, $false, $false, $false, $false, $true, 1, $false, "הנחה זה מאפשרת לנו להגדיר תהליך דטרמיניסטי של x_t-1 מ x_t באמצעות מודל שמאומן לשערך x_0 (הדאטה התחלתי מ-x_t). כלומר בכל איטרציה אנו קודם כל משערכים את x_0 באמצעות המודל ולאחר מכן בונים בצורה דטרמיניסטי אנו מחשבים x_t-1 מ-x_0 המשוערך. ", 2) | Out-Null

# 7. Rewrite the paragraph about speeding up generation
$d.Content.Find.Execute("המאמר המקורי על CM השתמש במרחק הנקרא LPIPS המודל דמיון סמנטי בין התמונות (דרך השוואה של אקטיבציות של מודלים מאומנים על דאטהסטים ענקיים של תמונות). המאמר המקורי גם התשמש ב-EMA (החלקה מעריכית) של משקלי המודל בתור המודל עבור x_t. יש כמובן חשיבות לבחירת השונות של האיטרציות.", $true, $false, $false, $false, $false, $true, 1, $false, "אבל איך זה בעצם כאשר לזירוז של תהליך יצירת הדאטה? מתברר ששערוך של x_t דרך שערוך x_0 מאפשר להקטין משמעותית את מספר האיטרציות וככה הדאטה נוצר מהם יותר.", 2) | Out-Null

# 8. Rewrite the closing paragraph
$d.Content.Find.Execute("אז המאמר שסוקרים היום משפר את תהליך האימון. השינוי הראשון הוא משקול של המרחקים כפונקציה של איטרציה t; ככל שמתקרבים ל 0 המשקול עולה. דבר שני זה שינוי של פונקציית מרחק מ-LPIPS לפונקציית הובר (Huber) עם טוויסט קטן. הדבר האחרון והמעניין הוא ביטול של EMA ל-x_t כלומר ההשוואה מתבצעת בין שני מודלים ״טהורים״ ל- x_t ו- x_t+1. גם הייפר פרמטרים אחרים עבור שינוי למשל השוניות של הרעש באיטרציות.", $true, $false, $false, $false, $false, $true, 1, $false, "מאמר מאוד מעניין - הסברתי אותו ממש בגדול, חובת קריאה לכל מי שאוהב מודלי דיפוזיה גנרטיביים.", 2) | Out-Null

# 9. Remove the "בקיצור..." paragraph entirely (it is dropped in the new version)
$d.Paragraphs(9).Range.Delete()

# 10. Update the arXiv link (now the last paragraph, after the deletion above)
$d.Content.Find.Execute("https://arxiv.org/abs/2310.14189", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/pdf/2010.02502", 2) | Out-Null
